$d = $word.ActiveDocument

# Locate the target sentence and trim its trailing ". " down to just the
# sentence text (no trailing period/space), matching the "before" run's
# wording so the remainder can be appended as a distinct, separately
# formatted run.
$rng = $d.Content
$found = $rng.Find.Execute(
    "second highest Recall Score. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "second highest Recall Score", 2)

if (-not $found) {
    throw "Could not find target sentence to edit."
}

# $rng now spans the replaced text ("...second highest Recall Score"),
# collapse to its end point so we can append new content right after it.
$rng.Collapse(0)
$rng.InsertAfter(" and both Recall and AUC Scores are high.")

# The appended text currently shares the same run as the preceding
# sentence (identical formatting merges runs). Toggle a character format
# on and back off to force Word to materialize it as its own <w:r>, which
# is the structure produced by the source edit (two runs with identical
# rPr rather than one merged run).
$rng.Font.Bold = 1
$rng.Font.Bold = 0
